$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that sits right after the
# title paragraph (it is being relocated to the end of the document). ---
$d.Paragraphs(2).Range.Delete() | Out-Null

# --- Step 2: insert a new paragraph right before the final "Prompt: ..." paragraph,
# containing a bold run with the page title. ---
$lastIdx = $d.Paragraphs.Count
$minBetPara = $d.Paragraphs($lastIdx - 1)

# Create a clean, empty paragraph right after the "Minimum bet..." bullet
# (and therefore right before the "Prompt:" paragraph), with no inherited styling.
$r = $minBetPara.Range
$r.Collapse(0)
$r.InsertAfter("`r") | Out-Null

# That new empty paragraph is now the second-to-last paragraph; fill it in with
# the same two-run shape ("<w:r/>" + bold text run) used elsewhere in the doc.
$newIdx = $lastIdx
$newPara = $d.Paragraphs($newIdx)
$npr = $newPara.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dolphin Gold with Stellar Jackpots Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$npr.InsertXML($xml) | Out-Null

# --- Step 3: replace the text of the (now last) "Prompt: ..." run with the
# relocated meta-description copy, keeping its italic formatting intact. ---
$d.Content.Find.Execute(
    "Prompt: Create a cartoon-style feature image for Dolphin Gold with Stellar Jackpots that features a happy Maya warrior wearing glasses. The image should have an underwater theme with the dolphin and gold elements incorporated into the background. The Maya warrior should be holding a treasure chest and smiling at the viewer. Use bright colors and bold lines to make the image pop and attract attention to the game's exciting features. The image should convey the idea of adventure and treasure while also showcasing the game's playful and enjoyable aspects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Dolphin Gold with Stellar Jackpots. Play for free and enjoy exciting gameplay, stunning graphics, and triple jackpots.",
    2
) | Out-Null

Write-Host "Edit complete."
